$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style/format from H1 (existing header cell) onto the new
# header cells I1 and J1, then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data column values for rows 2-4
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 5
